{"js": "// Update the 25 three-digit-by-one-digit multiplication prompts in the\n// practice-sheet table. Each old expression is unique in the document, so a\n// simple scoped search-and-replace (matchCase, whole text of the run) is\n// unambiguous for every cell.\nconst replacements = [\n  [\"857\u00d77=\", \"879\u00d79=\"],\n  [\"960\u00d77=\", \"204\u00d73=\"],\n  [\"462\u00d75=\", \"821\u00d76=\"],\n  [\"976\u00d77=\", \"336\u00d74=\"],\n  [\"987\u00d78=\", \"569\u00d77=\"],\n  [\"224\u00d73=\", \"772\u00d77=\"],\n  [\"624\u00d79=\", \"335\u00d74=\"],\n  [\"881\u00d78=\", \"673\u00d73=\"],\n  [\"743\u00d73=\", \"301\u00d74=\"],\n  [\"700\u00d73=\", \"830\u00d75=\"],\n  [\"114\u00d74=\", \"953\u00d79=\"],\n  [\"431\u00d74=\", \"804\u00d76=\"],\n  [\"952\u00d72=\", \"793\u00d72=\"],\n  [\"387\u00d73=\", \"640\u00d76=\"],\n  [\"349\u00d74=\", \"207\u00d79=\"],\n  [\"795\u00d75=\", \"759\u00d76=\"],\n  [\"484\u00d78=\", \"363\u00d74=\"],\n  [\"667\u00d79=\", \"657\u00d79=\"],\n  [\"428\u00d79=\", \"323\u00d73=\"],\n  [\"557\u00d78=\", \"404\u00d72=\"],\n  [\"707\u00d74=\", \"500\u00d76=\"],\n  [\"356\u00d72=\", \"625\u00d78=\"],\n  [\"475\u00d73=\", \"285\u00d73=\"],\n  [\"969\u00d77=\", \"121\u00d73=\"],\n  [\"893\u00d76=\", \"192\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 three-digit-by-one-digit multiplication prompts in the\n# practice-sheet table. Each old expression is unique in the document, so a\n# scoped Find/Replace (MatchCase on, whole document range) is unambiguous\n# for every cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"857\u00d77=\", \"879\u00d79=\"),\n    @(\"960\u00d77=\", \"204\u00d73=\"),\n    @(\"462\u00d75=\", \"821\u00d76=\"),\n    @(\"976\u00d77=\", \"336\u00d74=\"),\n    @(\"987\u00d78=\", \"569\u00d77=\"),\n    @(\"224\u00d73=\", \"772\u00d77=\"),\n    @(\"624\u00d79=\", \"335\u00d74=\"),\n    @(\"881\u00d78=\", \"673\u00d73=\"),\n    @(\"743\u00d73=\", \"301\u00d74=\"),\n    @(\"700\u00d73=\", \"830\u00d75=\"),\n    @(\"114\u00d74=\", \"953\u00d79=\"),\n    @(\"431\u00d74=\", \"804\u00d76=\"),\n    @(\"952\u00d72=\", \"793\u00d72=\"),\n    @(\"387\u00d73=\", \"640\u00d76=\"),\n    @(\"349\u00d74=\", \"207\u00d79=\"),\n    @(\"795\u00d75=\", \"759\u00d76=\"),\n    @(\"484\u00d78=\", \"363\u00d74=\"),\n    @(\"667\u00d79=\", \"657\u00d79=\"),\n    @(\"428\u00d79=\", \"323\u00d73=\"),\n    @(\"557\u00d78=\", \"404\u00d72=\"),\n    @(\"707\u00d74=\", \"500\u00d76=\"),\n    @(\"356\u00d72=\", \"625\u00d78=\"),\n    @(\"475\u00d73=\", \"285\u00d73=\"),\n    @(\"969\u00d77=\", \"121\u00d73=\"),\n    @(\"893\u00d76=\", \"192\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith,\n    # Replace(wdReplaceAll=2)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
